# Added email sending functionality - refresh contact list with new recipients
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$people = @(
    @("Timothy G Newton", " 312-349-8276", "monte.shanah@gmail.com", " Burr Ridge, Illinois(IL), 60527"),
    @("Lena T Chandler", " 559-324-9215", "delmer2013@yahoo.com", " Fresno, California(CA), 93710"),
    @("Rebecca T Sacks", " 630-722-9215", "sidney1979@gmail.com", " Oak Brook, Illinois(IL), 60523"),
    @("Bonita J Bagley", " 915-329-4474", "larry_krajc@gmail.com", " El Paso, Texas(TX), 79906"),
    @("Chester A Hicklin", " 585-243-9733", "sandrine2008@gmail.com", " Geneseo, New York(NY), 14454"),
    @("Christopher J Melendez", " 580-749-5751", "lindsey.kemm@hotmail.com", " Oklahoma City, Oklahoma(OK), 73129"),
    @("James M Dahle", " 915-241-5461", "alanna1975@yahoo.com", " El Paso, Texas(TX), 79912"),
    @("Sandra M Davis", " 310-239-0564", "orrin_gleichn@gmail.com", " Los Angeles, California(CA), 90017")
)

$row = 2
foreach ($person in $people) {
    $ws.Cells.Item($row, 1).Value = $person[0]
    $ws.Cells.Item($row, 2).Value = $person[1]
    $ws.Cells.Item($row, 3).Value = $person[2]
    $ws.Cells.Item($row, 4).Value = $person[3]
    $row++
}
